# Fire/Ice Bullet, Barrier & Turret Added
#
# This project-tracker sheet gets a new combined "防护罩、炮台" (Barrier & Turret)
# row inserted into the "玩家" (Player) section, just above "Wisplum精灵".
# The separate old "炮台" (Turret) and "杂项特效" (Misc effects) rows are removed
# (their text no longer appears anywhere in the sheet / shared strings table).
# A couple of existing rows (Ghost, 树桩状态机) also get a new "done/highlighted"
# green-on-green marker in column C, and several rows gain extra green-filled
# placeholder cells in columns C/D to line up with the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 20 (pushes "Wisplum精灵" / "火冰电元素表现" etc. down by one) ---
$ws.Rows(20).Insert()

# --- 2. The old "炮台" (now row 24) and "杂项特效" (now row 25) rows are fully
#        cleared (content + formatting) since their info is merged into the
#        new row 20 below; row 23 (old "防护罩") is likewise left empty. ---
$ws.Range("A23:F25").Clear()

# --- 3. Fill in the new combined row ---
$ws.Range("B20").Value = "防护罩、炮台"

# Green fill (style used elsewhere in the sheet as a "done" highlight) on C20/E20
$ws.Range("E9").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. New cells in the existing top rows (Ghost..Boss3 / 玩家 attacks) that
#        line up the C/D columns with the rest of the grid. ---
$ws.Range("E3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. "Done/highlighted" marker (green text on green fill) on C2 (Ghost) and
#        C28 (树桩状态机, shifted down by the row insert above). New font+fill
#        combo -> RGB(0,176,80) / #00B050, matching the sheet's existing green. ---
$ws.Range("C2").Interior.Color = 5287936
$ws.Range("C2").Font.Color = 5287936
$ws.Range("C28").Interior.Color = 5287936
$ws.Range("C28").Font.Color = 5287936

# --- 6. Restore the selection to where the author was last working ---
$ws.Range("J10").Select()
